# Update the cryptos list with freshly scraped values (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) are stored as plain text in this sheet (they can
# contain thousands separators like "42.156.85"). Whenever the new price
# looks like a plain number, force the cell to Text format first so Excel
# doesn't silently convert it to a numeric value and round-trip it with
# floating point noise (e.g. "305.66" -> 305.66000000000003) or drop
# trailing zeros (e.g. "2.90" -> 2.9).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "42.156.85"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.271.30"
$ws.Range("E3").Value = "  +0.71%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "305.66"
$ws.Range("E5").Value = "  +1.25%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "93.49"
$ws.Range("E6").Value = "  +1.33%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.16%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.489"
$ws.Range("E9").Value = "  +1.44%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "32.91"
$ws.Range("E10").Value = "  +1.82%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0805"
$ws.Range("E11").Value = "  +0.99%  "

# Rows 12-18 got reshuffled (ranking order changed) in addition to value
# updates - rewrite Coin/Link/Price/Volume for all seven rows.

# Row 12 - now TRON (was BinanceUSD)
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D12") "0.113"
$ws.Range("E12").Value = "  -1.68%  "

# Row 13 - now Polkadot (was TRON)
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "6.70"
$ws.Range("E13").Value = "  +0.63%  "

# Row 14 - now WrappedliquidstakedEther2.0 (was Polkadot)
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D14") "2.623.90"
$ws.Range("E14").Value = "  +0.78%  "

# Row 15 - now Chainlink (was WrappedliquidstakedEther2.0)
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D15") "14.33"
$ws.Range("E15").Value = "  +1.77%  "

# Row 16 - now WrappedEther (was Chainlink)
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "2.266.48"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17 - now Polygon (was WrappedEther)
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D17") "0.786"
$ws.Range("E17").Value = "  +3.90%  "

# Row 18 - now BinanceUSD (was Polygon)
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D18") "12.97"
$ws.Range("E18").Value = "  +1,196.22%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "42.013.09"
$ws.Range("E19").Value = "  +0.42%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D20") "12.72"
$ws.Range("E20").Value = "  +4.76%  "

# Row 21 - ShibaInu (price uses a subscript-3 digit to denote repeated zeros)
$sub3 = [char]0x2083
Set-TextValue $ws.Range("D21") "0.0${sub3}0919"
$ws.Range("E21").Value = "  +1.97%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.22%  "

# Row 23 - Litecoin
$ws.Range("E23").Value = "  +1.87%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "244.24"
$ws.Range("E24").Value = "  +1.41%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.60"
$ws.Range("E25").Value = "  +1.29%  "

# Row 26 - ImmutableX
Set-TextValue $ws.Range("D26") "1.94"
$ws.Range("E26").Value = "  +2.52%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.04%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.64%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +0.59%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -9.05%  "

# Row 31 - InjectiveProtocol
Set-TextValue $ws.Range("D31") "35.14"
$ws.Range("E31").Value = "  +3.58%  "

# Row 32 - Monero
Set-TextValue $ws.Range("D32") "160.19"
$ws.Range("E32").Value = "  +1.17%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.34"
$ws.Range("E33").Value = "  +3.65%  "

# Row 34 - FirstDigitalUSD
Set-TextValue $ws.Range("D34") "1.00"

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0744"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -0.68%  "

# Row 37 - Celestia
Set-TextValue $ws.Range("D37") "17.15"
$ws.Range("E37").Value = "  +4.17%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -1.03%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +1.44%  "

# Row 40 - Stellar
Set-TextValue $ws.Range("D40") "0.116"
$ws.Range("E40").Value = "  +0.46%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +0.17%  "

# Row 42 - RenderToken
Set-TextValue $ws.Range("D42") "4.00"
$ws.Range("E42").Value = "  +2.17%  "

# Row 43 - EnergySwap
$ws.Range("E43").Value = "  +0.44%  "

# Row 44 - Maker
Set-TextValue $ws.Range("D44") "2.014.49"
$ws.Range("E44").Value = "  -2.00%  "

# Row 45 - ApeXProtocol
$ws.Range("E45").Value = "  +9.44%  "

# Row 46 - VeChain
Set-TextValue $ws.Range("D46") "0.0284"
$ws.Range("E46").Value = "  +1.75%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  +2.10%  "

# Row 48 - NEARProtocol
Set-TextValue $ws.Range("D48") "2.90"
$ws.Range("E48").Value = "  +1.58%  "

# Row 49 - MultiversX
Set-TextValue $ws.Range("D49") "53.31"
$ws.Range("E49").Value = "  +3.39%  "

# Row 50 - Stacks
$ws.Range("E50").Value = "  +0.40%  "

# Row 51 - BitcoinSV
Set-TextValue $ws.Range("D51") "72.67"
$ws.Range("E51").Value = "  +2.92%  "
